$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9102510809898376
$ws.Range("B1").Value = 1.713655829429626
$ws.Range("C1").Value = 3.994026184082031
$ws.Range("D1").Value = 3.813498497009277
$ws.Range("E1").Value = 0.9171925783157349
